$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen" ---
$resumen = $wb.Worksheets.Item("Resumen")
$resumen.Range("B2").Value = "Z3"
$resumen.Range("C2").Value = 541.3548816088227

# --- Sheet "Solucion" (shuffle of the Salida/B column) ---
$solucion = $wb.Worksheets.Item("Solucion")
$solucion.Range("B5").Value = "S051"
$solucion.Range("B6").Value = "S031"
$solucion.Range("B9").Value = "S022"
$solucion.Range("B10").Value = "S035"
$solucion.Range("B12").Value = "S032"
$solucion.Range("B13").Value = "S053"
$solucion.Range("B17").Value = "S033"
$solucion.Range("B18").Value = "S026"
$solucion.Range("B19").Value = "S012"
$solucion.Range("B23").Value = "S014"
$solucion.Range("B26").Value = "S045"
$solucion.Range("B28").Value = "S025"
$solucion.Range("B29").Value = "S055"
$solucion.Range("B30").Value = "S015"
$solucion.Range("B31").Value = "S002"
$solucion.Range("B32").Value = "S046"
$solucion.Range("B34").Value = "S013"
$solucion.Range("B35").Value = "S056"
$solucion.Range("B38").Value = "S057"
$solucion.Range("B40").Value = "S036"
$solucion.Range("B41").Value = "S027"
$solucion.Range("B45").Value = "S018"
$solucion.Range("B46").Value = "S058"
$solucion.Range("B48").Value = "S008"
$solucion.Range("B53").Value = "S059"
$solucion.Range("B54").Value = "S019"
$solucion.Range("B56").Value = "S050"
$solucion.Range("B58").Value = "S030"
$solucion.Range("B59").Value = "S020"
$solucion.Range("B60").Value = "S060"
$solucion.Range("B61").Value = "S040"

# --- Sheet "Metricas" ---
$metricas = $wb.Worksheets.Item("Metricas")
$metricas.Range("B2").Value = 541.2268461455292
$metricas.Range("B3").Value = 541.3041247702455
$metricas.Range("B4").Value = 541.3548816088227
